$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "243.43"

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "23.80"

# Row 4
$ws.Cells.Item(4, 2).Value = "LEO"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "3.561"
$ws.Cells.Item(4, 5).Value = "3LEOLEO"

# Row 5
$ws.Cells.Item(5, 2).Value = "HuobiToken"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "5.261"
$ws.Cells.Item(5, 5).Value = "4HuobiTokenHT"

# Row 6
$ws.Cells.Item(6, 2).Value = "Cronos"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.05811"
$ws.Cells.Item(6, 5).Value = "5CronosCRO"

# Row 7
$ws.Cells.Item(7, 2).Value = "KuCoinToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "6.488"
$ws.Cells.Item(7, 5).Value = "6KuCoinTokenKCS"

# Row 8
$ws.Cells.Item(8, 2).Value = "GateToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.345"
$ws.Cells.Item(8, 5).Value = "7GateTokenGT"

# Row 9
$ws.Cells.Item(9, 2).Value = "MXToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.8081"
$ws.Cells.Item(9, 5).Value = "8MXTokenMX"

# Row 10
$ws.Cells.Item(10, 2).Value = "FTXToken"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8753"
$ws.Cells.Item(10, 5).Value = "9FTXTokenFTT"

# Row 11
$ws.Cells.Item(11, 2).Value = "One"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.01034"
$ws.Cells.Item(11, 5).Value = "10OneONEBestin24h"

# Row 12
$ws.Cells.Item(12, 2).Value = "WazirX"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.1389"
$ws.Cells.Item(12, 5).Value = "11WazirXWRX"

# Row 13
$ws.Cells.Item(13, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.07256"
$ws.Cells.Item(13, 5).Value = "12MandalaExchangeTokenMDX"

# Row 14
$ws.Cells.Item(14, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.03075"
$ws.Cells.Item(14, 5).Value = "13LiechtensteinCryptoassetsExchangeLCX"

# Row 15
$ws.Cells.Item(15, 2).Value = "BitrueCoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.03055"
$ws.Cells.Item(15, 5).Value = "14BitrueCoinBTR"

# Row 16
$ws.Cells.Item(16, 2).Value = "BitMartToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.09327"
$ws.Cells.Item(16, 5).Value = "15BitMartTokenBMX"

# Row 17
$ws.Cells.Item(17, 2).Value = "MCDex"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.857"
$ws.Cells.Item(17, 5).Value = "16MCDexMCB"

# Row 18
$ws.Cells.Item(18, 2).Value = "BitForexToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.001543"
$ws.Cells.Item(18, 5).Value = "17BitForexTokenBF"

# Row 19
$ws.Cells.Item(19, 2).Value = "CoinExToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.04701"
$ws.Cells.Item(19, 5).Value = "18CoinExTokenCET"

# Row 20
$ws.Cells.Item(20, 2).Value = "TigerCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.006169"
$ws.Cells.Item(20, 5).Value = "19TigerCashTCH"

# Row 21
$ws.Cells.Item(21, 2).Value = "BitKan"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.001266"
$ws.Cells.Item(21, 5).Value = "20BitKanKAN"

# Row 22
$ws.Cells.Item(22, 2).Value = "HotbitToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.004594"
$ws.Cells.Item(22, 5).Value = "21HotbitTokenHTB"

# Row 23
$ws.Cells.Item(23, 2).Value = "NitroEx"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.00008698"
$ws.Cells.Item(23, 5).Value = "22NitroExNTX"

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.165"

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.3174"

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1310"

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.0002343"

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.03793"

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.006319"

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1051"

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.002899"
$ws.Cells.Item(43, 5).Value = "42CEJICEJIWorstin24h"

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.007981"

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00005535"

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5999"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOIN"

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.01329"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.00002099"

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0001999"
